$d = $word.ActiveDocument

# The "Java options" section used to list two reference hyperlinks and then
# a "Splint" heading followed by an "Email from Splint people" sub-heading.
# Per the commit message ("Removed extra headings that I'm not going to go
# into"), remove:
#   - the "http://pmd.sourceforge.net/" hyperlink paragraph
#   - the "http://checkstyle.sourceforge.net/index.html" hyperlink paragraph
#   - the blank paragraph that immediately follows those two links
#   - the "Email from Splint people" Heading 3 paragraph
# (the "Splint" Heading 2 paragraph itself is kept).
#
# Paragraphs are located by their text content (rather than hard-coded
# indices) and collected first; the actual deletions happen afterwards,
# from the highest index down to the lowest, so that removing one
# paragraph never invalidates the index of another one still pending.

$toDelete = New-Object System.Collections.ArrayList
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text

    if ($text -like "*http://pmd.sourceforge.net/*") {
        [void]$toDelete.Add($i)
    }
    elseif ($text -like "*http://checkstyle.sourceforge.net/index.html*") {
        [void]$toDelete.Add($i)
        if (($i + 1) -le $count) {
            # the blank paragraph right after the checkstyle link goes too
            [void]$toDelete.Add($i + 1)
        }
    }
    elseif ($text -like "*Email from Splint people*") {
        [void]$toDelete.Add($i)
    }
}

$unique = $toDelete | Sort-Object -Unique
$ordered = @($unique)
for ($j = $ordered.Count - 1; $j -ge 0; $j--) {
    $idx = $ordered[$j]
    $d.Paragraphs.Item($idx).Range.Delete()
}

Write-Host "Removed $($ordered.Count) paragraphs."
